# Mark traits as "currently in deck" by setting n_of_cards (column G) to 1
# for every trait row that doesn't already track a card count, and tidy up
# leftover formatting on the MOL sheet's now-unused G/H columns.

$wb = $excel.ActiveWorkbook

$traits = $wb.Worksheets.Item("traits")
$mol = $wb.Worksheets.Item("MOL")

# --- traits sheet: set G2:G170 to 1, skipping rows that already carry a count ---
$traits.Range("G2:G39").Value = 1
$traits.Range("G41:G104").Value = 1
$traits.Range("G106:G170").Value = 1

# --- MOL sheet: clear the stray border formatting left on G3:H7 ---
$mol.Range("G3:H7").Clear()

# --- selections / active sheet, to match the end-of-session UI state ---
$mol.Range("G12").Select()

$traits.Activate()
$traits.Range("A107").Select()
$traits.Range("G121:G170").Select()
$traits.Cells.Item(121, 7).Activate()
